$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VaxEfficacyRandVal")

$values = @(
    0.98302901785187102,
    0.97569909610059602,
    0.96661691542948802,
    0.96449336079088199,
    0.94966027691702903,
    0.992714189094005,
    0.98527778542775402,
    0.98073908405315802,
    0.95276479738130904,
    0.96517960579049,
    0.98294133220660895,
    0.94764711149346204,
    0.99738600015754397,
    0.93488189171523906,
    0.98130292320127099,
    0.99063105173537902,
    0.97538098617963898,
    0.99013459154161498,
    0.98596822114529803,
    0.984211409284106,
    0.99595648393467695,
    0.92443543793937799,
    0.96358490499643801,
    0.991421136988483,
    0.95146410767954204,
    0.97732868521638905,
    0.91950310395177803,
    0.92576811046260399,
    0.93577149011205196,
    0.94442478115643802,
    0.97095597095934105,
    0.97759225803972905,
    0.97932299820070501,
    0.98540887483963302,
    0.97511592048051998,
    0.99271738719040603,
    0.88542299210890696,
    0.99125685264181795,
    0.97981986307005098,
    0.98117666706284601,
    0.929072850608333,
    0.971605876422038,
    0.95955799089417204,
    0.96991711011649095,
    0.97036738475080297,
    0.97558465983076004,
    0.94257075293368897,
    0.96469515473795597,
    0.96111428112824204,
    0.95380710297596405,
    0.98190038150445702,
    0.99004425828125098,
    0.99067710564327405,
    0.95736639110752497,
    0.950333690683559,
    0.96722472063138298,
    0.98167390663607001,
    0.99855105206241201,
    0.90894520982402005,
    0.98564824152854802,
    0.98529713699147903,
    0.97707688557274197,
    0.99228114739644901,
    0.96250221069644704,
    0.97322424858703105,
    0.96881156591151896,
    0.95353125645708703,
    0.97910212946813302,
    0.93876611578412095,
    0.97544677312374894,
    0.95607964595454098,
    0.95101425511876103,
    0.97806836671944897,
    0.91018555309710103,
    0.94475132778253401,
    0.96017277083785402,
    0.98426225576308402,
    0.98376612032337096,
    0.97584182913752804,
    0.99013963618517298,
    0.95940626150571495,
    0.97612361454225505,
    0.97857339647901398,
    0.98920005060975102,
    0.97043463694231502,
    0.98684651508293297,
    0.99584228821217502,
    0.97931751948377199,
    0.94111563409874399,
    0.99552497984493904,
    0.970760946882406,
    0.95747277435417699,
    0.965362009635921,
    0.98448760502085997,
    0.98780602538000095,
    0.99011184644424699,
    0.96492430236343396,
    0.98914598646981999,
    0.98166435253159601,
    0.97636222212081503
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $values[$i]
}

$ws.Range("F19").Select() | Out-Null